$wb = $excel.ActiveWorkbook

# --- "white meat" sheet ---
# healthy flag becomes a real boolean TRUE (was string "yes")
$wsWhiteMeat = $wb.Worksheets.Item("white meat")
$wsWhiteMeat.Range("B3").Value = $true

# --- "red meat" sheet ---
# healthy flag becomes a real boolean FALSE (was string "no")
$wsRedMeat = $wb.Worksheets.Item("red meat")
$wsRedMeat.Range("B3").Value = $false

# add a new "healthy substitute" -> "plant-based red meat" row
$wsRedMeat.Range("A7").Value = "healthy substitute"
$wsRedMeat.Range("B7").Value = "plant-based red meat"

# --- update selections (previous cursor positions) on the relevant sheets ---
# fish, fowl and meat sheets keep their existing selection, only these move:
$wsPork = $wb.Worksheets.Item("pork")
$wsPork.Range("A4").Select()

$wsRedMeat.Range("D19").Select()

$wsWhiteMeat.Range("C1").Select()
